# #5: property aircraft done
#
# The legislator's aircraft (航空器) property sheet is removed entirely
# (it had no real holdings), and two rows elsewhere that had been
# mistakenly tagged with the "land" property_category are corrected to
# their real category ("building" for the 建物 sheet, "car" for the
# 汽車 sheet).

$wb = $excel.ActiveWorkbook

# Remove the whole "航空器" (aircraft) worksheet.
$wsAircraft = $wb.Worksheets.Item("航空器")
$wsAircraft.Delete()

# Fix the property_category column (I) on the 建物 (building) sheet --
# it was incorrectly carrying "land" for both data rows.
$wsBuilding = $wb.Worksheets.Item("建物")
$wsBuilding.Range("I2").Value = "building"
$wsBuilding.Range("I3").Value = "building"

# Fix the property_category column (H) on the 汽車 (car) sheet --
# it was incorrectly carrying "land" for both data rows.
$wsCar = $wb.Worksheets.Item("汽車")
$wsCar.Range("H2").Value = "car"
$wsCar.Range("H3").Value = "car"
